# Tilty Quad BT sheet: remove the obsolete BOM line "1 uF Capacitor / C5, C9, C12"
# (row 26). Deleting the whole row shifts every row below it up by one and
# Excel re-points the SUM()/shared-formula ranges and cross-sheet shared-string
# references automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tilty Quad BT")
$ws.Rows(26).Delete()

# Restore the selection shown in the saved file after the edit.
$null = $ws.Range("K31").Select()
